# Fixed Typo in BOM - Line 21 Trough -> Through
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in B21 ("SI4825 Trough Hole KIT PCB" -> "Si4825 Through Hole Kit PCB")
$ws.Range("B21").Value = "Si4825 Through Hole Kit PCB"

# Clear the stray "?" placeholder values in F17/G17
$ws.Range("F17").ClearContents()
$ws.Range("G17").ClearContents()

# Move the active selection to C22 (as recorded in the saved file)
$ws.Range("C22").Select()
